# Registeration.xlsx edit script
# - Replace garbled "password" column values (C2:C6) with literal "password"
#   on both the "Login" and "Regx" sheets.
# - Change active sheet/selection: "Login" becomes the active (tab-selected) sheet
#   with selection D11; "Regx" is no longer tab-selected, selection D10.

$wb = $excel.ActiveWorkbook

# --- Sheet "Login" ---
$loginSheet = $wb.Worksheets.Item("Login")
$loginSheet.Range("C2:C6").Value = "password"

# --- Sheet "Regx" ---
$regxSheet = $wb.Worksheets.Item("Regx")
$regxSheet.Range("C2:C6").Value = "password"

# --- Selections / active sheet ---
$regxSheet.Range("D10").Select()
$loginSheet.Activate()
$loginSheet.Range("D11").Select()
